$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Kota Lahir"), pushing old D (Tanggal Lahir) -> E
# and old E (password) -> F.
$ws.Range("D1").EntireColumn.Insert()

# --- Row 1 (headers) ---
$ws.Range("D1").Value = "Kota Lahir"
$ws.Range("E1").Value = "Tanggal Lahir"
$ws.Range("F1").Value = "password"

# --- Row 2 ---
$ws.Range("A2").Value = "Kemal S"
$ws.Range("B2").Value = 2341760196
$ws.Range("C2").Value = "MAHASISWA"
$ws.Range("D2").Value = "Surabaya"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = 37632
$ws.Range("E2").NumberFormat = "m/d/yyyy"
$ws.Range("F2").Value = 2341760196
$ws.Range("F2").ClearFormats()

# --- Row 3 ---
$ws.Range("A3").Value = "Ismi Atika"
$ws.Range("B3").Value = 2341760036
$ws.Range("C3").Value = "MAHASISWA"
$ws.Range("D3").Value = "Bangkalan"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "15-01-2004"
$ws.Range("E3").NumberFormat = "m/d/yyyy"
$ws.Range("F3").Value = 2341760036

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 16.5703125
$ws.Columns.Item(4).ColumnWidth = 13.42578125
$ws.Columns.Item(5).ColumnWidth = 15.5703125
$ws.Columns.Item(6).ColumnWidth = 16.85546875

# --- Selection ---
$ws.Range("E3").Select()
